# Update Sheets via scheduled runner: refresh currentAveragePrice / LevePrice /
# LeveProfit figures across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H64").Value = 4011.28
$ws.Range("I64").Value = 3949.2
$ws.Range("J64").Value = 4259.6
$ws.Range("K64").Value = 3949.2
$ws.Range("L64").Value = 4259.6
$ws.Range("M64").Value = -3701.2
$ws.Range("N64").Value = -4755.6
$ws.Range("H67").Value = 4011.28
$ws.Range("I67").Value = 3949.2
$ws.Range("J67").Value = 4259.6
$ws.Range("K67").Value = 3949.2
$ws.Range("L67").Value = 4259.6
$ws.Range("M67").Value = -3091.2
$ws.Range("N67").Value = -5975.6
$ws.Range("H113").Value = 2906.7646
$ws.Range("I113").Value = 2616.875
$ws.Range("K113").Value = 2616.875
$ws.Range("M113").Value = 637.125
$ws.Range("H132").Value = 3335543.2
$ws.Range("I132").Value = 2025.0428
$ws.Range("K132").Value = 6075.1284
$ws.Range("M132").Value = -3545.1284
$ws.Range("H137").Value = 883.8
$ws.Range("I137").Value = 847.1053000000001
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 2541.3159
$ws.Range("L137").Value = 3000
$ws.Range("M137").Value = 8.684099999999944
$ws.Range("N137").Value = -8100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1375
$ws.Range("I2").Value = 1100
$ws.Range("J2").Value = 2200
$ws.Range("K2").Value = 1100
$ws.Range("L2").Value = 2200
$ws.Range("M2").Value = -987
$ws.Range("N2").Value = -2426
$ws.Range("H74").Value = 1362.9474
$ws.Range("I74").Value = 830.72
$ws.Range("J74").Value = 2386.4614
$ws.Range("K74").Value = 830.72
$ws.Range("L74").Value = 2386.4614
$ws.Range("M74").Value = 43.27999999999997
$ws.Range("N74").Value = -4134.4614
$ws.Range("H77").Value = 1362.9474
$ws.Range("I77").Value = 830.72
$ws.Range("J77").Value = 2386.4614
$ws.Range("K77").Value = 4153.6
$ws.Range("L77").Value = 11932.307
$ws.Range("M77").Value = 214.3999999999996
$ws.Range("N77").Value = -20668.307
$ws.Range("H97").Value = 871.375
$ws.Range("I97").Value = 493.33334
$ws.Range("J97").Value = 2005.5
$ws.Range("K97").Value = 493.33334
$ws.Range("L97").Value = 2005.5
$ws.Range("M97").Value = 2.666659999999979
$ws.Range("N97").Value = -2997.5
$ws.Range("H102").Value = 1818.125
$ws.Range("I102").Value = 1557.8334
$ws.Range("K102").Value = 1557.8334
$ws.Range("M102").Value = 64.16660000000002
$ws.Range("H116").Value = 1375
$ws.Range("I116").Value = 1100
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 1100
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = 1194
$ws.Range("N116").Value = -6788
$ws.Range("H132").Value = 2192.3333
$ws.Range("I132").Value = 1688.1515
$ws.Range("J132").Value = 4965.3335
$ws.Range("K132").Value = 5064.4545
$ws.Range("L132").Value = 14896.0005
$ws.Range("M132").Value = -2534.4545
$ws.Range("N132").Value = -19956.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1375
$ws.Range("I3").Value = 1100
$ws.Range("J3").Value = 2200
$ws.Range("K3").Value = 1100
$ws.Range("L3").Value = 2200
$ws.Range("M3").Value = -986
$ws.Range("N3").Value = -2428
$ws.Range("H94").Value = 797.7568
$ws.Range("I94").Value = 633.1786
$ws.Range("J94").Value = 1309.7778
$ws.Range("K94").Value = 633.1786
$ws.Range("L94").Value = 1309.7778
$ws.Range("M94").Value = -182.1786
$ws.Range("N94").Value = -2211.7778
$ws.Range("H105").Value = 3312.0852
$ws.Range("I105").Value = 1638.1818
$ws.Range("K105").Value = 1638.1818
$ws.Range("M105").Value = 108.8181999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3587.4644
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3587.4644
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3587.4644
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -4177.4644
$ws.Range("H34").Value = 3587.4644
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 3587.4644
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3587.4644
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -3991.4644

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 868189.5
$ws.Range("J2").Value = 1984311
$ws.Range("L2").Value = 11905866
$ws.Range("N2").Value = -11906092
$ws.Range("H107").Value = 804.4857
$ws.Range("J107").Value = 948.3570999999999
$ws.Range("L107").Value = 2845.0713
$ws.Range("N107").Value = -6685.0713
$ws.Range("H132").Value = 1552.7059
$ws.Range("I132").Value = 700.8889
$ws.Range("J132").Value = 2511
$ws.Range("K132").Value = 6308.0001
$ws.Range("L132").Value = 22599
$ws.Range("M132").Value = -3778.0001
$ws.Range("N132").Value = -27659

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 14850.692
$ws.Range("J135").Value = 14850.692
$ws.Range("L135").Value = 14850.692
$ws.Range("N135").Value = -24990.692
$ws.Range("H140").Value = 49808.777
$ws.Range("J140").Value = 49808.777
$ws.Range("L140").Value = 49808.777
$ws.Range("N140").Value = -60168.777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5792
$ws.Range("I7").Value = 2048.6428
$ws.Range("J7").Value = 9823.308000000001
$ws.Range("K7").Value = 2048.6428
$ws.Range("L7").Value = 9823.308000000001
$ws.Range("M7").Value = -1936.6428
$ws.Range("N7").Value = -10047.308
$ws.Range("H126").Value = 5792
$ws.Range("I126").Value = 2048.6428
$ws.Range("J126").Value = 9823.308000000001
$ws.Range("K126").Value = 6145.928400000001
$ws.Range("L126").Value = 29469.924
$ws.Range("M126").Value = -3675.928400000001
$ws.Range("N126").Value = -34409.924

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 16226.596
$ws.Range("J64").Value = 16226.596
$ws.Range("L64").Value = 16226.596
$ws.Range("N64").Value = -16722.596
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 16226.596
$ws.Range("J67").Value = 16226.596
$ws.Range("L67").Value = 16226.596
$ws.Range("N67").Value = -17942.596
$ws.Range("H126").Value = 1439.8704
$ws.Range("I126").Value = 1371.8889
$ws.Range("J126").Value = 1575.8334
$ws.Range("K126").Value = 4115.6667
$ws.Range("L126").Value = 4727.5002
$ws.Range("M126").Value = -1645.6667
$ws.Range("N126").Value = -9667.5002
$ws.Range("H136").Value = 1252.5156
$ws.Range("I136").Value = 647.4231
$ws.Range("J136").Value = 3874.5833
$ws.Range("K136").Value = 1942.2693
$ws.Range("L136").Value = 11623.7499
$ws.Range("M136").Value = 607.7307000000001
$ws.Range("N136").Value = -16723.7499
